$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'65.441.93"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = "'2.938.50"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.58%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = "'568.58"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.81%  '
$ws.Range('D6').Value = "'158.79"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.87%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('E8').Value = '  -0.55%  '
$ws.Range('D9').Value = "'2.935.31"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.50%  '
$ws.Range('D10').Value = "'6.73"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.32%  '
$ws.Range('E11').Value = '  -3.42%  '
$ws.Range('D12').Value = "'0.460"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.84%  '
$ws.Range('D13').Value = "'0.0000245"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.42%  '
$ws.Range('D14').Value = "'34.36"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.72%  '
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').Value = "'65.420.36"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('D17').Value = "'3.425.00"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.61%  '
$ws.Range('D18').Value = "'7.00"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('D19').Value = "'2.937.56"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.62%  '
$ws.Range('D20').Value = "'15.71"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +13.14%  '
$ws.Range('D21').Value = "'444.72"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.86%  '
$ws.Range('D22').Value = "'0.692"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.93%  '
$ws.Range('D23').Value = "'7.29"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.03%  '
$ws.Range('D24').Value = "'82.25"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').Value = "'2.25"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('D26').Value = "'12.13"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.86%  '
$ws.Range('D27').Value = "'10.06"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.01%  '
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').Value = "'2.36"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.38%  '
$ws.Range('E31').Value = '  -1.38%  '
$ws.Range('D32').Value = "'0.0000101"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.97%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = "'0.111"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.23%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = "'27.09"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.30%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = "'0.972"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.29%  '
$ws.Range('D37').Value = "'5.75"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.19%  '
$ws.Range('D38').Value = "'49.69"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('E39').Value = '  +3.09%  '
$ws.Range('E40').Value = '  -8.71%  '
$ws.Range('D41').Value = "'0.301"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.78%  '
$ws.Range('E42').Value = '  -2.05%  '
$ws.Range('E43').Value = '  -7.29%  '
$ws.Range('D44').Value = "'8.49"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('D45').Value = "'383.04"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.32%  '
$ws.Range('E46').Value = '  -0.82%  '
$ws.Range('D47').Value = "'2.698.45"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.67%  '
$ws.Range('D48').Value = "'133.44"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('D50').Value = "'2.20"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.96%  '
$ws.Range('E51').Value = '  +0.34%  '
